$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 12:22"

# Row 11 - Estados Unidos
$ws.Range("B11").Value = 76389
$ws.Range("C11").Value = 1512
$ws.Range("D11").Value = 49933
$ws.Range("E11").Value = 21679
$ws.Range("F11").Value = 3643
$ws.Range("G11").Value = 94
$ws.Range("H11").Value = 4777

# Row 20 - only F changes
$ws.Range("F20").Value = 232

# Row 31 - Rumania (label unchanged)
$ws.Range("B31").Value = 7216
$ws.Range("C31").Value = 337
$ws.Range("D31").Value = 1217
$ws.Range("E31").Value = 5637
$ws.Range("F31").Value = 245

# Row 32 - label swaps from Dinamarca to Noruega, with new data values
$ws.Range("A32").Value = "Noruega"
$ws.Range("B32").Value = 6686
$ws.Range("C32").Value = 63
$ws.Range("D32").Value = 32
$ws.Range("E32").Value = 6512
$ws.Range("F32").Value = 59
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 142

# Row 33 - label swaps from Noruega to Dinamarca, with new data values
$ws.Range("A33").Value = "Dinamarca"
$ws.Range("B33").Value = 6681
$ws.Range("C33").Value = 170
$ws.Range("D33").Value = 2515
$ws.Range("E33").Value = 3867
$ws.Range("F33").Value = 100
$ws.Range("H33").Value = 299

# Row 90
$ws.Range("B90").Value = 658
$ws.Range("C90").Value = 17
$ws.Range("E90").Value = 556
